$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.2669677734375
$ws.Range("B1").Value = 2.744543313980103
$ws.Range("C1").Value = 5.03374719619751
$ws.Range("D1").Value = 2.00143575668335
$ws.Range("E1").Value = 1.032090425491333
